$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 229
$ws.Range("I33").Value = 97.833336
$ws.Range("J33").Value = 386.4
$ws.Range("K33").Value = 97.833336
$ws.Range("L33").Value = 386.4
$ws.Range("M33").Value = 131.166664
$ws.Range("N33").Value = -844.4
$ws.Range("H62").Value = 326501.06
$ws.Range("J62").Value = 4426.25
$ws.Range("L62").Value = 4426.25
$ws.Range("N62").Value = -5674.25
$ws.Range("H65").Value = 326501.06
$ws.Range("J65").Value = 4426.25
$ws.Range("L65").Value = 22131.25
$ws.Range("N65").Value = -28371.25
$ws.Range("H98").Value = 1856.6364
$ws.Range("I98").Value = 1925.7142
$ws.Range("J98").Value = 406
$ws.Range("K98").Value = 1925.7142
$ws.Range("L98").Value = 406
$ws.Range("M98").Value = -427.7141999999999
$ws.Range("N98").Value = -3402
$ws.Range("H112").Value = 1810.7142
$ws.Range("I112").Value = 766.6667
$ws.Range("J112").Value = 2095.4546
$ws.Range("K112").Value = 2300.0001
$ws.Range("L112").Value = 6286.3638
$ws.Range("M112").Value = -1192.0001
$ws.Range("N112").Value = -8502.363799999999
$ws.Range("H116").Value = 4116.8184
$ws.Range("I116").Value = 3945.2942
$ws.Range("J116").Value = 4700
$ws.Range("K116").Value = 3945.2942
$ws.Range("L116").Value = 4700
$ws.Range("M116").Value = -503.2941999999998
$ws.Range("N116").Value = -11584
$ws.Range("H122").Value = 1856.6364
$ws.Range("I122").Value = 1925.7142
$ws.Range("J122").Value = 406
$ws.Range("K122").Value = 5777.142599999999
$ws.Range("L122").Value = 1218
$ws.Range("M122").Value = -3327.142599999999
$ws.Range("N122").Value = -6118
$ws.Range("H129").Value = 1078.5264
$ws.Range("J129").Value = 1230.6451
$ws.Range("L129").Value = 3691.9353
$ws.Range("N129").Value = -13691.9353
$ws.Range("H137").Value = 1037.579
$ws.Range("I137").Value = 978.6923
$ws.Range("J137").Value = 1165.1666
$ws.Range("K137").Value = 2936.0769
$ws.Range("L137").Value = 3495.4998
$ws.Range("M137").Value = -386.0769
$ws.Range("N137").Value = -8595.4998
$ws.Range("H138").Value = 2981.732
$ws.Range("I138").Value = 1045.3889
$ws.Range("J138").Value = 4124.4917
$ws.Range("K138").Value = 3136.1667
$ws.Range("L138").Value = 12373.4751
$ws.Range("M138").Value = 2003.8333
$ws.Range("N138").Value = -22653.4751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2176.56
$ws.Range("I2").Value = 1810.4445
$ws.Range("J2").Value = 3118
$ws.Range("K2").Value = 1810.4445
$ws.Range("L2").Value = 3118
$ws.Range("M2").Value = -1697.4445
$ws.Range("N2").Value = -3344
$ws.Range("H32").Value = 2346
$ws.Range("I32").Value = 2346
$ws.Range("K32").Value = 2346
$ws.Range("M32").Value = -2059
$ws.Range("H61").Value = 2082.913
$ws.Range("I61").Value = 2116.4736
$ws.Range("J61").Value = 1923.5
$ws.Range("K61").Value = 2116.4736
$ws.Range("L61").Value = 1923.5
$ws.Range("M61").Value = -1904.4736
$ws.Range("N61").Value = -2347.5
$ws.Range("H74").Value = 1288.0869
$ws.Range("I74").Value = 1301.1818
$ws.Range("K74").Value = 1301.1818
$ws.Range("M74").Value = -427.1818000000001
$ws.Range("H77").Value = 1288.0869
$ws.Range("I77").Value = 1301.1818
$ws.Range("K77").Value = 6505.909000000001
$ws.Range("M77").Value = -2137.909000000001
$ws.Range("H88").Value = 3023.6667
$ws.Range("I88").Value = 2649.5
$ws.Range("J88").Value = 3130.5715
$ws.Range("K88").Value = 2649.5
$ws.Range("L88").Value = 3130.5715
$ws.Range("M88").Value = -2243.5
$ws.Range("N88").Value = -3942.5715
$ws.Range("H91").Value = 3023.6667
$ws.Range("I91").Value = 2649.5
$ws.Range("J91").Value = 3130.5715
$ws.Range("K91").Value = 2649.5
$ws.Range("L91").Value = 3130.5715
$ws.Range("M91").Value = -1245.5
$ws.Range("N91").Value = -5938.5715
$ws.Range("H116").Value = 2176.56
$ws.Range("I116").Value = 1810.4445
$ws.Range("J116").Value = 3118
$ws.Range("K116").Value = 1810.4445
$ws.Range("L116").Value = 3118
$ws.Range("M116").Value = 483.5554999999999
$ws.Range("N116").Value = -7706
$ws.Range("H117").Value = 23248
$ws.Range("J117").Value = 23248
$ws.Range("L117").Value = 23248
$ws.Range("N117").Value = -32426
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 44473.332
$ws.Range("J120").Value = 44473.332
$ws.Range("L120").Value = 44473.332
$ws.Range("N120").Value = -54149.332
$ws.Range("H121").Value = 13112
$ws.Range("J121").Value = 13112
$ws.Range("L121").Value = 13112
$ws.Range("N121").Value = -16606
$ws.Range("H122").Value = 1091.25
$ws.Range("I122").Value = 843.8889
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 2531.6667
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -81.66670000000022
$ws.Range("N122").Value = -10400.0002
$ws.Range("H123").Value = 49980
$ws.Range("J123").Value = 49980
$ws.Range("L123").Value = 49980
$ws.Range("N123").Value = -59780
$ws.Range("H132").Value = 1427.6364
$ws.Range("I132").Value = 1273.6154
$ws.Range("J132").Value = 1999.7142
$ws.Range("K132").Value = 3820.8462
$ws.Range("L132").Value = 5999.142599999999
$ws.Range("M132").Value = -1290.8462
$ws.Range("N132").Value = -11059.1426
$ws.Range("H136").Value = 2082.913
$ws.Range("I136").Value = 2116.4736
$ws.Range("J136").Value = 1923.5
$ws.Range("K136").Value = 6349.4208
$ws.Range("L136").Value = 5770.5
$ws.Range("M136").Value = -3799.4208
$ws.Range("N136").Value = -10870.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2176.56
$ws.Range("I3").Value = 1810.4445
$ws.Range("J3").Value = 3118
$ws.Range("K3").Value = 1810.4445
$ws.Range("L3").Value = 3118
$ws.Range("M3").Value = -1696.4445
$ws.Range("N3").Value = -3346
$ws.Range("H105").Value = 3895.8276
$ws.Range("I105").Value = 3927.4285
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 3927.4285
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -2180.4285
$ws.Range("N105").Value = -6505
$ws.Range("H134").Value = 1409.4694
$ws.Range("I134").Value = 1230.3112
$ws.Range("J134").Value = 3425
$ws.Range("K134").Value = 3690.9336
$ws.Range("L134").Value = 10275
$ws.Range("M134").Value = -1155.9336
$ws.Range("N134").Value = -15345

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 67328.375
$ws.Range("I31").Value = 6072
$ws.Range("J31").Value = 169422.33
$ws.Range("K31").Value = 6072
$ws.Range("L31").Value = 169422.33
$ws.Range("M31").Value = -5777
$ws.Range("N31").Value = -170012.33
$ws.Range("H34").Value = 67328.375
$ws.Range("I34").Value = 6072
$ws.Range("J34").Value = 169422.33
$ws.Range("K34").Value = 6072
$ws.Range("L34").Value = 169422.33
$ws.Range("M34").Value = -5870
$ws.Range("N34").Value = -169826.33
$ws.Range("H58").Value = 2357.0571
$ws.Range("I58").Value = 829.913
$ws.Range("J58").Value = 5284.0835
$ws.Range("K58").Value = 829.913
$ws.Range("L58").Value = 5284.0835
$ws.Range("M58").Value = -626.913
$ws.Range("N58").Value = -5690.0835
$ws.Range("H62").Value = 4217.1113
$ws.Range("I62").Value = 2630.8
$ws.Range("J62").Value = 6200
$ws.Range("K62").Value = 2630.8
$ws.Range("L62").Value = 6200
$ws.Range("M62").Value = -2006.8
$ws.Range("N62").Value = -7448
$ws.Range("H65").Value = 4217.1113
$ws.Range("I65").Value = 2630.8
$ws.Range("J65").Value = 6200
$ws.Range("K65").Value = 13154
$ws.Range("L65").Value = 31000
$ws.Range("M65").Value = -10034
$ws.Range("N65").Value = -37240
$ws.Range("H122").Value = 1282.625
$ws.Range("I122").Value = 1243.5
$ws.Range("K122").Value = 3730.5
$ws.Range("M122").Value = -1280.5
$ws.Range("H136").Value = 2357.0571
$ws.Range("I136").Value = 829.913
$ws.Range("J136").Value = 5284.0835
$ws.Range("K136").Value = 2489.739
$ws.Range("L136").Value = 15852.2505
$ws.Range("M136").Value = 60.26099999999997
$ws.Range("N136").Value = -20952.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 3375
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3375
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10125
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -10579
$ws.Range("H34").Value = 1127.5333
$ws.Range("I34").Value = 358.25
$ws.Range("J34").Value = 1407.2727
$ws.Range("K34").Value = 1074.75
$ws.Range("L34").Value = 4221.8181
$ws.Range("M34").Value = -990.75
$ws.Range("N34").Value = -4389.8181
$ws.Range("H39").Value = 10992453
$ws.Range("J39").Value = 10992453
$ws.Range("L39").Value = 32977359
$ws.Range("N39").Value = -32977947
$ws.Range("H132").Value = 889.25
$ws.Range("I132").Value = 873.4286
$ws.Range("K132").Value = 7860.8574
$ws.Range("M132").Value = -5330.8574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1669.45
$ws.Range("I16").Value = 1827.6428
$ws.Range("J16").Value = 1300.3334
$ws.Range("K16").Value = 1827.6428
$ws.Range("L16").Value = 1300.3334
$ws.Range("M16").Value = -1657.6428
$ws.Range("N16").Value = -1640.3334
$ws.Range("H46").Value = 876.6875
$ws.Range("J46").Value = 629.6667
$ws.Range("L46").Value = 629.6667
$ws.Range("N46").Value = -1005.6667
$ws.Range("H132").Value = 3008.1667
$ws.Range("I132").Value = 2410
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 7230
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -4700
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 305.4737
$ws.Range("J113").Value = 385
$ws.Range("L113").Value = 1155
$ws.Range("N113").Value = -5495
$ws.Range("H132").Value = 2424
$ws.Range("I132").Value = 2151
$ws.Range("J132").Value = 2860.8
$ws.Range("K132").Value = 6453
$ws.Range("L132").Value = 8582.400000000001
$ws.Range("M132").Value = -3923
$ws.Range("N132").Value = -13642.4
$ws.Range("H136").Value = 1350.3928
$ws.Range("I136").Value = 1522.4348
$ws.Range("J136").Value = 559
$ws.Range("K136").Value = 4567.3044
$ws.Range("L136").Value = 1677
$ws.Range("M136").Value = -2017.3044
$ws.Range("N136").Value = -6777
